$wb = $excel.ActiveWorkbook

# --- "Flags" sheet: update the Categories / AllColors flag values ---
$flags = $wb.Worksheets.Item("Flags")
$flags.Range("B3").Value = "Debug"
$flags.Range("B4").Value = "False"

# --- "Tests" sheet: clear out the now-unused test row (row 42) ---
$tests = $wb.Worksheets.Item("Tests")
$tests.Range("B42:D42").ClearContents()
